# The workbook originally has an unused leading label column (A) with
# row numbers (1, 3, 10, 18) duplicated into column F; that leading
# column is being removed so the data shifts left by one column
# (old B->A, C->B, D->C, E->D, F->E), and the "MODEL_CONDITION" header
# text is being corrected to "MODELCONDITION".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading column A entirely; Excel shifts everything else left.
$ws.Columns.Item(1).Delete()

# Fix the header text that now lives in D1.
$ws.Range("D1").Value2 = "MODELCONDITION"
